$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing the existing rows 25-26 down to 26-27
$ws.Rows.Item(25).Insert()

# Fill the newly inserted row 25 with this week's data (same shape as the row below,
# new date and new volume)
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 45209
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100107
$ws.Range("H25").Value = "Otros"
$ws.Range("I25").Value = 100107002
$ws.Range("J25").Value = "Chirimoya"
$ws.Range("K25").Value = "Cultivar IV Región"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 22000
$ws.Range("O25").Value = 22000
$ws.Range("P25").Value = 22000
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("R25").Value = "Provincia de Limarí"
$ws.Range("S25").Value = 2200
$ws.Range("T25").Value = 10

# Preserve the date number format on the newly inserted date cell (same style as D26/D27)
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
